# Insert a new row at row 470 (pushes existing rows 470-546 down to 471-547)
# and populate it with the new data record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(470).Insert()

$ws.Cells.Item(470, 1).Value = 5
$ws.Cells.Item(470, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(470, 3).Value = "Maule"
$ws.Cells.Item(470, 4).Value = 45180
$ws.Cells.Item(470, 5).Value = 7
$ws.Cells.Item(470, 6).Value = "Fruta"
$ws.Cells.Item(470, 7).Value = 100102
$ws.Cells.Item(470, 8).Value = "Cítricos"
$ws.Cells.Item(470, 9).Value = 100102004
$ws.Cells.Item(470, 10).Value = "Mandarina"
$ws.Cells.Item(470, 11).Value = "Murcott"
$ws.Cells.Item(470, 12).Value = "Primera"
$ws.Cells.Item(470, 13).Value = 250
$ws.Cells.Item(470, 14).Value = 8000
$ws.Cells.Item(470, 15).Value = 8000
$ws.Cells.Item(470, 16).Value = 8000
$ws.Cells.Item(470, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(470, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(470, 19).Value = 444
$ws.Cells.Item(470, 20).Value = 18
